# Change timezone to Pakistan time: update Nouman's Out Time ("D3") to a
# new time value, and bump the In Time ("C3") forward accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "8:12:55 PM"
$ws.Range("D3").Value = "8:12:56 PM"
